$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "a MSC"
$ws.Range("B2").Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/behaviour-change-msc"
$ws.Range("D2").Value = "A postgraduate degree in behaviour change centred around the systematic application of behaviour change theory and methods to design, implement and evaluate interventions, primarily using the Behaviour Change Wheel. This approach equips students to work in this emerging and exciting field to address social, health and environmental challenges."

$ws.Range("A3").Value = "B msc"
$ws.Range("B3").Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/architectural-design"
$ws.Range("D3").Value = "Architectural Design at The Bartlett is invested in the frontiers of advanced architecture and design and its convergence with science and technology. Composed of an international staff of experts and students, this programme is designed to deliver diverse yet focused strands of speculative research, emphasising the key role computation plays within complex design synthesis."
